# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") rows 2-16
$newK = @{
    2  = 5
    3  = 6
    4  = 3
    5  = 2
    6  = 9
    7  = 2
    8  = 3
    9  = 1
    10 = 5
    11 = 3
    12 = 3
    13 = 4
    14 = 4
    15 = 4
    16 = 3
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
